$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 774, shifting existing rows 774:815 down to 775:816.
$ws.Rows.Item(774).EntireRow.Insert()

# Populate the newly inserted row 774 with the new data point.
# Force column A to be treated as plain text (avoid Excel's automatic
# date-string -> serial-number conversion), then clear the formatting
# that the text-coercion trick leaves behind so the cell ends up with
# no explicit style, matching the rest of the column.
$ws.Range("A774").NumberFormat = "@"
$ws.Range("A774").Value = "2026/02/03"
$ws.Range("A774").ClearFormats()

$ws.Range("B774").Value = "火"
$ws.Range("C774").Value = 16
$ws.Range("D774").Value = 201

Write-Host "row inserted and populated"
